# Remove the grid of "thumb_*" montage picture shapes from slide 1,
# leaving only the title, content placeholder, and the QR-code picture
# (Picture 3) untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -like "Picture *" -and $shape.Name -ne "Picture 3") {
        $shape.Delete()
    }
}
